$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header suffixes: _old -> _FV2210, _new -> _FV2304
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace '_old$', '_FV2210')
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace '_new$', '_FV2304')
}

# Freeze the header row (split below row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# Turn the data range into an Excel Table ("Table1")
$range = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
